# Insert a new weekly price row at row 30 (pushing the existing rows 30-56
# down to 31-57) and populate it with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 30, shifting rows 30-56 -> 31-57.
$ws.Rows.Item(30).Insert()

# Fill the newly inserted row 30 with the new data point.
$ws.Cells.Item(30, 1).Value = 11
$ws.Cells.Item(30, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(30, 3).Value = "Bíobío"
$ws.Cells.Item(30, 4).Value = 45118
$ws.Cells.Item(30, 5).Value = 8
$ws.Cells.Item(30, 6).Value = 100114007
$ws.Cells.Item(30, 7).Value = "Jengibre"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 30
$ws.Cells.Item(30, 11).Value = 15000
$ws.Cells.Item(30, 12).Value = 15000
$ws.Cells.Item(30, 13).Value = 15000
$ws.Cells.Item(30, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(30, 15).Value = "Perú"
$ws.Cells.Item(30, 16).Value = 1154
$ws.Cells.Item(30, 17).Value = 13
$ws.Cells.Item(30, 18).Value = "Hortaliza"
